$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model-name order for column A (rows 2-26) after the training refresh
$modelNames = @{
    2 = "model_14_7_0"
    3 = "model_14_7_22"
    4 = "model_14_7_21"
    5 = "model_14_7_20"
    6 = "model_14_7_19"
    7 = "model_14_7_18"
    8 = "model_14_7_17"
    9 = "model_14_7_16"
    10 = "model_14_7_15"
    11 = "model_14_7_14"
    12 = "model_14_7_13"
    13 = "model_14_7_23"
    14 = "model_14_7_12"
    15 = "model_14_7_10"
    16 = "model_14_7_9"
    17 = "model_14_7_8"
    18 = "model_14_7_7"
    19 = "model_14_7_6"
    20 = "model_14_7_5"
    21 = "model_14_7_4"
    22 = "model_14_7_3"
    23 = "model_14_7_2"
    24 = "model_14_7_1"
    25 = "model_14_7_11"
    26 = "model_14_7_24"
}

# Metric values (B..Q) shared by every data row, as produced by the refreshed
# ensemble training run. [double] casts keep exponent literals well-formed.
$metricCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
$metricVals = @([double]"0.999997025781996", [double]"0.9990399586905852", [double]"0.9999999999999789", [double]"0.9999945721615463", [double]"0.9999988831789071", [double]"2.776301953652792e-06", [double]"0.0008961564213956752", [double]"2.521337463430942e-14", [double]"1.910840072249828e-06", [double]"9.554200487316011e-07", [double]"0.0001004477823987822", [double]"0.001666223860606009", [double]"1.000007931248011", [double]"0.001737158459910895", [double]"91.58878149661076", [double]"131.8116837172614")

foreach ($r in $modelNames.Keys) {
    $ws.Range("A$r").Value = $modelNames[$r]
    for ($i = 0; $i -lt $metricCols.Length; $i++) {
        $ws.Range("$($metricCols[$i])$r").Value = $metricVals[$i]
    }
}
